# edit.ps1 - apply the "wizard" template update to the epicrisis docx.
#
# Two kinds of edits are needed:
#  1. Two run-merges in the "Находился на лечении" paragraph:
#       "{patient.arrivalAt" + " | date:"   -> "{patient.arrivalAt | date:"
#       "{patient.departureAt" + " | date:" -> "{patient.departureAt | date:"
#     plus a double-space -> single-space fix just before the first merge.
#  2. Eighteen insertions of one or more new runs holding a human readable
#     Russian label immediately before a lone "{placeholder}" run (each such
#     run is the only run in its paragraph).
#
# Plain Find.Execute(..., Replace:=wdReplaceOne) normalizes/merges whole
# runs of identically-formatted text around the edit point in this runtime,
# which would over-merge neighbouring runs that the diff wants to keep
# separate. InsertXML on a *collapsed* Range does not have that problem -
# inserting a "<w:p>...</w:p>" fragment at the exact start of an existing
# paragraph splices its runs in front of the existing content in place
# (no new paragraph mark is produced when the insertion point is already a
# paragraph boundary), so that is used throughout.

$d = $word.ActiveDocument

$wordNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rPrEn = '<w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>'
$rPrRu = '<w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:rtl w:val="0"/><w:lang w:val="ru-RU"/></w:rPr>'

function New-Run($rPr, $text, $preserve) {
    if ($preserve) {
        return '<w:r>' + $rPr + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        return '<w:r>' + $rPr + '<w:t>' + $text + '</w:t></w:r>'
    }
}

# Insert one or more runs (already-built "<w:r>...</w:r>" XML, concatenated)
# immediately before the unique-in-document $anchorText, which must be the
# first (and, for our 18 simple cases, only) run of its paragraph.
function Insert-RunsBeforeText($doc, $anchorText, $runsXml) {
    $rng = $doc.Content
    $ok = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Insert-RunsBeforeText: anchor not found: $anchorText"
    }
    $insRng = $doc.Range($rng.Start, $rng.Start)
    $xml = '<w:p ' + $wordNs + '>' + $runsXml + '</w:p>'
    $insRng.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Fix the "Находился на лечении..." paragraph: tidy the double space
#    and merge the two run-pairs that build the arrivalAt/departureAt
#    template expressions.
# ---------------------------------------------------------------------

$rng = $d.Content
$ok = $rng.Find.Execute("Находился", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "anchor paragraph 'Находился' not found"
}
$rng.Expand(4) | Out-Null
$insPos = $rng.Start
$delRng = $d.Range($rng.Start, $rng.End - 1)
$delRng.Delete()

$quote = [char]0x2019

$parts = @()
$parts += New-Run $rPrRu "Находился  на лечении в " $true
$parts += New-Run $rPrRu "1 -" $false
$parts += New-Run $rPrRu "ом кардиологическом отделении с " $true
$parts += New-Run $rPrEn "{patient.arrivalAt | date:" $false
$parts += New-Run $rPrEn $quote $false
$parts += New-Run $rPrEn "D.MM.Y" $false
$parts += New-Run $rPrEn $quote $false
$parts += New-Run $rPrEn "}" $false
$parts += New-Run $rPrRu "г" $false
$parts += New-Run $rPrRu ". " $true
$parts += New-Run $rPrRu "по " $true
$parts += New-Run $rPrEn "{patient.departureAt | date:" $false
$parts += New-Run $rPrEn $quote $false
$parts += New-Run $rPrEn "D.MM.Y" $false
$parts += New-Run $rPrEn $quote $false
$parts += New-Run $rPrEn "}" $false
$parts += New-Run $rPrRu "г" $false
$parts += New-Run $rPrRu "." $false

$paraXml = '<w:p ' + $wordNs + '>' + ([string]::Join("", $parts)) + '</w:p>'
$insRng = $d.Range($insPos, $insPos)
$insRng.InsertXML($paraXml)

# ---------------------------------------------------------------------
# 2) Insert descriptive-label runs in front of the (lone) placeholder run
#    of eighteen single-run paragraphs.
# ---------------------------------------------------------------------

Insert-RunsBeforeText $d "{usdScopia.obp}" (New-Run $rPrEn "УЗИ ОБП " $true)

Insert-RunsBeforeText $d "{usdScopia.bca}" (New-Run $rPrEn "УЗИ БЦА" $false)

Insert-RunsBeforeText $d "{usdScopia.fgds}" (New-Run $rPrEn "ФГДС " $true)

Insert-RunsBeforeText $d "{usdScopia.bronchoscopy}" (New-Run $rPrEn "Бронхоскопия " $true)

Insert-RunsBeforeText $d "{xray.crgogk}" (New-Run $rPrEn "КРГ ОГК " $true)

Insert-RunsBeforeText $d "{xray.crgSkull}" (New-Run $rPrEn "КРГ черепа в двух проекциях " $true)

$joints = @()
$joints += New-Run $rPrEn "Рентгенография правого" $false
$joints += New-Run $rPrEn "/" $false
$joints += New-Run $rPrEn "левого коленного" $false
$joints += New-Run $rPrEn "/" $false
$joints += New-Run $rPrEn "тазобедренного сустава " $true
Insert-RunsBeforeText $d "{xray.jointsRoentgenography}" ([string]::Join("", $joints))

Insert-RunsBeforeText $d "{ct.head}" (New-Run $rPrEn "РКТ головы " $true)

Insert-RunsBeforeText $d "{ct.ogk}" (New-Run $rPrEn "РКТ ОГК " $true)

$ctobp = @()
$ctobp += New-Run $rPrEn "РКТ ОБП с контрастирование" $false
$ctobp += New-Run $rPrEn "/" $false
$ctobp += New-Run $rPrEn "без контрастирования " $true
Insert-RunsBeforeText $d "{ct.obp}" ([string]::Join("", $ctobp))

Insert-RunsBeforeText $d "{examination.ophthalmologist}" (New-Run $rPrEn "Осмотр офтальмолога " $true)

$ent = @()
$ent += New-Run $rPrEn "Осмотр лор" $false
$ent += New-Run $rPrEn "-" $false
$ent += New-Run $rPrEn "врача " $true
Insert-RunsBeforeText $d "{examination.entDoctor}" ([string]::Join("", $ent))

Insert-RunsBeforeText $d "{examination.urological}" (New-Run $rPrEn "Осмотр уролога " $true)

Insert-RunsBeforeText $d "{examination.physiotherapist}" (New-Run $rPrEn "Осмотр физеотерапевта " $true)

Insert-RunsBeforeText $d "{examination.psychiatric}" (New-Run $rPrEn "Осмотр психиатра " $true)

Insert-RunsBeforeText $d "{examination.surgeon}" (New-Run $rPrEn "Осмотр хирурга " $true)

Insert-RunsBeforeText $d "{examination.oncologist}" (New-Run $rPrEn "Осмотр онколога " $true)

Insert-RunsBeforeText $d "{examination.midwife}" (New-Run $rPrEn "Осмотр акушерки " $true)

Write-Output "edit.ps1 completed"
